$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates: C2 and E2 change from text-stored numbers to real numbers,
# F2's stock count changes from 1 to 5.
$ws.Range("C2").Value = 15
$ws.Range("E2").Value = 13500
$ws.Range("F2").Value = 5

# New row 3 data (NUTRIBON / PERRO item). C3 and E3 stay text-stored numbers,
# matching the same "number stored as text" pattern used elsewhere in the sheet.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "NUTRIBON"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "20"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "PERRO"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8300"
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").Value = 14
